# 19-04-21: Design Complete. First Order Placed and Payed for
#
# Re-creates the save-time state captured in the target diff for the
# "Bill of Materials-Communication" sheet:
#   - row 1's height was shrunk from the tall wrapped-header height
#     (56.25pt) down to a compact 20.25pt custom height
#   - the last active selection on the sheet was moved from K48 to C2
#     (i.e. the user clicked back on cell C2 before saving/closing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink header row (row 1) to a custom height of 20.25 points.
$ws.Rows.Item(1).RowHeight = 20.25

# Move the saved selection/active cell to C2.
[void]$ws.Range("C2").Select()
